$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet holds an "Estado de Cuenta" (account statement) table of workers
# in arrears. The edit removes the previous set of rows and replaces them
# with an updated worker/period ledger (7 workers, 9 periods, 11 data rows
# instead of 8), and refreshes the totals that summarize that table.
# ---------------------------------------------------------------------------

# --- Grow the data table from 8 rows (16-23) to 11 rows (16-26): insert 3
# rows just above the last ("closing border") row, copying the format of
# the row right above the insertion point so the new rows pick up the same
# (non-bordered) table-body style instead of a blank default style.
$ws.Range("B23:J25").Insert(-4121)  # xlShiftDown
$ws.Range("B22:J22").Copy()
$ws.Range("B23:J25").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Replace the data rows 16-26 with the new worker/period ledger ---
function Set-LedgerRow($r, $doc, $name, $period, $mora, $salario) {
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = $doc
    $ws.Cells.Item($r, 4).Value = $name
    $ws.Cells.Item($r, 5).Value = $period
    $ws.Cells.Item($r, 6).Value = $mora
    $ws.Cells.Item($r, 7).Value = $salario
}

Set-LedgerRow 16 "9097544"    "JOHN HAROL SEPULVEDA ZABALETA" "2103" 52000 1300000
Set-LedgerRow 17 "9097544"    "JOHN HAROL SEPULVEDA ZABALETA" "2102" 38133 1300000
Set-LedgerRow 18 "73118049"   "CARLOS ENRIQUE RESTREPO GOMEZ" "2505" 1898  1423500
Set-LedgerRow 19 "1101813420" "DAGOBERTO MUNZON OSORIO"       "1702" 27578 800000
Set-LedgerRow 20 "92255665"   "OCTAVIO AUGUSTO MONTES HOYOS"  "2101" 48000 1200000
Set-LedgerRow 21 "92255665"   "OCTAVIO AUGUSTO MONTES HOYOS"  "1704" 29509 1200000
Set-LedgerRow 22 "8641226"    "FABIAN ALBERTO SALAS CORONADO" "1908" 9600  1200000
Set-LedgerRow 23 "104300053"  "CARLOS ENRIQUE PELAEZ AVILA"   "1910" 40000 1000000
Set-LedgerRow 24 "104300053"  "CARLOS ENRIQUE PELAEZ AVILA"   "1909" 40000 1000000
Set-LedgerRow 25 "104300053"  "CARLOS ENRIQUE PELAEZ AVILA"   "1908" 5333  1000000
Set-LedgerRow 26 "1143404454" "JOSE MANUEL MEJIA BOLIVAR"     "2101" 35112 877803

# --- Refresh the summary fields above the table ---
$ws.Range("E11").Value = 327163   # VALOR MORA total
$ws.Range("C13").Value = 7        # Cant. Trabajadores
$ws.Range("F13").Value = 9        # Cant. Periodos

# --- Column D widened slightly to fit the new (longer) name strings ---
$ws.Columns("D").ColumnWidth = 33.8
